$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so the stale shared strings ("yolima", "ggg",
# "rr", "r", "e") and the stale cell layout (old 2-row A1:F2 block) are
# dropped, matching the new A1:H4 data block.
$ws.Cells.Clear()

# --- Row 1 ---
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "hola"
$ws.Range("C1").Value = "hola mundo"
$ws.Range("D1").Value = "hola"
$ws.Range("F1").Value = "hola"

# --- Row 2 ---
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "mundo"
$ws.Range("C2").Value = "hola mundo"
$ws.Range("D2").Value = "nn"
$ws.Range("F2").Value = "hola"

# --- Row 3 ---
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "pc"
$ws.Range("E3").Value = "g"

# --- Row 4 ---
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "mary"
$ws.Range("E4").Value = "*"

# Extra record field added afterwards, re-using the "hola"/"hola mundo"
# strings already interned above, plus two new ones.
$ws.Range("C3").Value = "ii"

# A numeric-looking value that must be kept as literal text, so it is
# forced onto a text format before assignment (same trick Excel's COM
# automation uses: NumberFormat "@" prevents the "0.999" -> number
# auto-conversion that a plain .Value assignment would trigger).
$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "0.999"

# Genuine numeric cells.
$ws.Range("H1").Value = 0.78
$ws.Range("C4").Value = 67
$ws.Range("D4").Value = 67

$ws.Range("K9").Select()
